# Apply cryptos.xlsx price/volume update (commit: "Updated cryptos list on Sat May  6 09:56:52 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.384.87'
$ws.Range("E2").Value = '  +0.79%  '

$ws.Range("D3").Value = '1.942.88'
$ws.Range("E3").Value = '  +2.07%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.20'
$ws.Range("E5").Value = '  -0.12%  '

$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4626'
$ws.Range("E7").Value = '  +0.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3872'
$ws.Range("E8").Value = '  -0.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.98'
$ws.Range("E9").Value = '  -0.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07826'
$ws.Range("E10").Value = '  -0.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9746'
$ws.Range("E11").Value = '  -1.49%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.59'
$ws.Range("E12").Value = '  +3.09%  '

$ws.Range("D13").Value = '1.949.65'
$ws.Range("E13").Value = '  +3.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.084'
$ws.Range("E14").Value = '  +0.51%  '

$ws.Range("E15").Value = '  -0.08%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07025'
$ws.Range("E16").Value = '  +0.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.64'
$ws.Range("E17").Value = '  -1.45%  '

$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009813'
$ws.Range("E19").Value = '  -1.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.08'
$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("D22").Value = '29.403.07'
$ws.Range("E22").Value = '  +0.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.483'
$ws.Range("E23").Value = '  +3.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.05'
$ws.Range("E24").Value = '  -0.68%  '

$ws.Range("D25").Value = '2.172.22'
$ws.Range("E25").Value = '  +2.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.096'
$ws.Range("E26").Value = '  -0.33%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.13'
$ws.Range("E27").Value = '  +0.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.39'
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.742'
$ws.Range("E29").Value = '  -2.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.46'
$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.857'
$ws.Range("E31").Value = '  -0.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09357'
$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8596'
$ws.Range("E33").Value = '  -3.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.172'
$ws.Range("E34").Value = '  -1.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.306'
$ws.Range("E35").Value = '  -1.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.102'
$ws.Range("E36").Value = '  -1.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05765'
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("E38").Value = '  -1.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02081'
$ws.Range("E39").Value = '  -0.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.678'
$ws.Range("E40").Value = '  +0.13%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5672'
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1779'
$ws.Range("E42").Value = '  -0.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.397'
$ws.Range("E43").Value = '  -3.23%  '

$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000002848'
$ws.Range("E44").Value = '  +35.59%  '

$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.728'
$ws.Range("E45").Value = '  +6.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5289'
$ws.Range("E46").Value = '  -1.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.42'
$ws.Range("E47").Value = '  -3.74%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06871'
$ws.Range("E48").Value = '  -2.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.083'
$ws.Range("E49").Value = '  -5.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.816'
$ws.Range("E50").Value = '  -1.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.24'
$ws.Range("E51").Value = '  -1.67%  '

